$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the full contents of row 2 and row 3 (including which
# optional cells are populated), as if the two species/observation
# records had traded places in the sheet.

# --- Row 2 gets the values that used to be in row 3 ---
$ws.Range("A2").Value = 86851042
$ws.Range("B2").Value = 96251
$ws.Range("E2").Value = 219790
$ws.Range("F2").Value = "Fläcknycklar"
$ws.Range("G2").Value = "Dactylorhiza maculata"
$ws.Range("H2").Value = "(L.) Soó"
$ws.Range("P2").Value = "Längs E4 mellan Kongberget och Gnarp, Hls"
$ws.Range("Q2").Value = 615689.084506036
$ws.Range("R2").Value = 6862637.86594828
$ws.Range("S2").Value = 10
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2018-06-27"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2018-06-27"
$ws.Range("AS2").Value = ""
$ws.Range("AW2").Value = "Sofia Lundman"
$ws.Range("AX2").Value = "Sofia Lundman, Oskar Wallströmer"

# --- Row 3 gets the values that used to be in row 2 ---
$ws.Range("A3").Value = 105312389
$ws.Range("B3").Value = 78503
$ws.Range("E3").Value = 6456
$ws.Range("F3").Value = "Skinnlav"
$ws.Range("G3").Value = "Leptogium saturninum"
$ws.Range("H3").Value = "(Dicks.) Nyl."
$ws.Range("P3").Value = "Storrönningen, Hls"
$ws.Range("Q3").Value = 616060.6447056353
$ws.Range("R3").Value = 6863194.20718522
$ws.Range("S3").Value = 5
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2022-06-03"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2022-06-03"
$ws.Range("AS3").Value = "Henrik Tykosson"
$ws.Range("AW3").Value = "Helene Andersson"
$ws.Range("AX3").Value = "Henrik Tykosson"
